$wb = $excel.ActiveWorkbook
$th = $wb.Theme()
$cs = $th.ThemeColorScheme()
$item3 = $cs.Colors(3)
$item3.RGB = 0x7d491f
$v = $item3.RGB()
Write-Output "rgb readback = $v"
